# Add the new "ML_AUC" worksheet (as the last tab) and populate it with the
# GLM vs. SMOTE model-performance summary table, matching the authored
# "Add files via upload" commit.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new worksheet after the last existing sheet -----------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "ML_AUC"

# --- 2. Fill in the data -----------------------------------------------------
# Populate the cells in the same order the shared-string table records them
# in (column headers first, then the row labels top-to-bottom) so the
# resulting sharedStrings.xml ordering matches the authored workbook.
$ws.Range("C1").Value = "SMOTE"
$ws.Range("B1").Value = "GLM"

$ws.Range("A3").Value = "AUC"
$ws.Range("A4").Value = "Accuracy"
$ws.Range("A5").Value = "Specificity"
$ws.Range("A6").Value = "Sensitvity"
$ws.Range("A2").Value = "Threshold"

$ws.Range("B2").Value = 0.0093903323604451003
$ws.Range("C2").Value = 0.0079950401805517698

$ws.Range("B3").Value = 0.73150000000000004
$ws.Range("C3").Value = 0.73480000000000001

$ws.Range("B4").Value = 0.89706091154337597
$ws.Range("C4").Value = 0.874769274456908

$ws.Range("B5").Value = 0.90086083213773305
$ws.Range("C5").Value = 0.87804878040000001

$ws.Range("B6").Value = 0.534246575342466
$ws.Range("C6").Value = 0.56164383561643805

# --- 3. Restore each sheet's own cursor/selection ----------------------------
# (".Select()" also activates its sheet, so set the non-final sheets first
# and finish on ML_AUC so it ends up the active tab, matching the workbook's
# saved activeTab.)
$wsTraining = $wb.Worksheets.Item("ARM_Training")
$wsTraining.Range("B18").Select() | Out-Null

$wsModel = $wb.Worksheets.Item("ML_Model")
$wsModel.Range("L13").Select() | Out-Null

$ws.Range("E7").Select() | Out-Null
